$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove column H: unmerge the G:H merged ranges first, then clear column H ---
$ws.Range("A1:H1").UnMerge()
$ws.Range("G2:H2").UnMerge()
$ws.Range("G4:H4").UnMerge()
$ws.Range("G5:H5").UnMerge()
$ws.Range("G6:H6").UnMerge()
$ws.Range("G7:H7").UnMerge()

$ws.Columns.Item(8).Clear()

# Re-merge the title row across the now-smaller range A1:G1
$ws.Range("A1:G1").Merge()

# --- Row 1 (title) becomes left-aligned (was center/default horizontal before) ---
$ws.Range("A1:G1").HorizontalAlignment = -4131   # xlLeft

# --- G5 no longer carries a formula/value (third position cell is blank) ---
$ws.Range("G5").Formula = ""

# --- G4 / G6 formula-result cells lose their right alignment ---
$ws.Range("G4").HorizontalAlignment = 1          # xlGeneral
$ws.Range("G6").HorizontalAlignment = 1          # xlGeneral

# --- G7 (sum) also loses its right alignment ---
$ws.Range("G7").HorizontalAlignment = 1          # xlGeneral

# --- Update the active selection shown in the sheet view ---
$ws.Range("H11").Select()
